$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.014.83'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.265.47'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.18'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.59'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.44%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.600'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.88%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.836.88'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.44'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '68.050.66'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.07%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.226.13'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.72'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.29'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '416.36'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +7.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.53'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.12%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.45'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.68'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.81%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.64%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.86%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '162.72'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.44'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.88'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.00'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.796'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.35'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.639.29'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.80'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0675'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.42'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '336.98'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.33'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.14%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.27'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.976'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.32%  '
